$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from serial 45212 (2023-10-13) to serial 45221 (2023-10-22),
# keeping the existing date formatting/style on the cells.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
